# Add hotspot support for 8.X+ to xlsx (fix #133)
#
# The "Metrics" sheet (last tab, driven by the "metrics" Excel Table on
# range A1:A2 with a single "Colonne1" header) is duplicated to create a
# new "Hotspots" tab positioned immediately before "Metrics". The
# duplicate keeps the original "metrics" table/data (renamed to
# "hotspots"), while the original sheet becomes the new "Metrics" tab and
# gets its own fresh "metrics" table so both tabs keep working
# independently (same shape: a one-column autofiltered table used as a
# named range for later exports).

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

# Duplicate "Metrics" -> places the copy right after the original, named
# "Metrics (2)" for now. The ORIGINAL sheet/table stay intact and become
# "Hotspots"; the NEW copy becomes the fresh "Metrics" tab.
$metrics.Copy($null, $metrics)

$copy = $wb.Worksheets.Item("Metrics (2)")

# The original sheet (with its original "metrics" table) turns into
# "Hotspots".
$metrics.Name = "Hotspots"
$hotspotsTable = $metrics.ListObjects.Item("metrics")
$hotspotsTable.Name = "hotspots"

# The duplicated sheet becomes the new "Metrics" tab; rebuild its own
# "metrics" table (the Copy() didn't bring the ListObject along) over the
# same A1:A2 range, matching the other sheets' look and feel.
$copy.Name = "Metrics"
$newTable = $copy.ListObjects.Add(1, $copy.Range("A1:A2"), $null, 1)
$newTable.Name = "metrics"
$newTable.TableStyle = "TableStyleLight16"

# Leave the workbook focused on the new "Hotspots" tab (mirrors the
# previously-active "Metrics" tab position).
$metrics.Activate()
$metrics.Select()
$metrics.Range("A2").Select()
